# Refresh cryptos list values (price & 1h volume change) to the latest snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = '59.466.80'
$cell.ClearFormats()
$ws.Cells.Item(2, 5).Value = '  +1.01%  '

# Row 3
$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.530.88'
$cell.ClearFormats()
$ws.Cells.Item(3, 5).Value = '  +0.34%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.09%  '

# Row 5
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '536.91'
$cell.ClearFormats()
$ws.Cells.Item(5, 5).Value = '  +0.58%  '

# Row 6
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = '141.78'
$cell.ClearFormats()
$ws.Cells.Item(6, 5).Value = '  -1.24%  '

# Row 7
$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.995'
$cell.ClearFormats()
$ws.Cells.Item(7, 5).Value = '  -0.37%  '

# Row 8
$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.567'
$cell.ClearFormats()
$ws.Cells.Item(8, 5).Value = '  -0.66%  '

# Row 9
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.540.73'
$cell.ClearFormats()
$ws.Cells.Item(9, 5).Value = '  +0.79%  '

# Row 10
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0997'
$cell.ClearFormats()
$ws.Cells.Item(10, 5).Value = '  +0.07%  '

# Row 11
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.161'
$cell.ClearFormats()
$ws.Cells.Item(11, 5).Value = '  +1.80%  '

# Row 12
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.42'
$cell.ClearFormats()
$ws.Cells.Item(12, 5).Value = '  -3.04%  '

# Row 13
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.356'
$cell.ClearFormats()
$ws.Cells.Item(13, 5).Value = '  +1.82%  '

# Row 14
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.983.78'
$cell.ClearFormats()
$ws.Cells.Item(14, 5).Value = '  +0.81%  '

# Row 15
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = '23.43'
$cell.ClearFormats()
$ws.Cells.Item(15, 5).Value = '  -1.62%  '

# Row 16
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = '59.482.61'
$cell.ClearFormats()
$ws.Cells.Item(16, 5).Value = '  +1.09%  '

# Row 17
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0000141'
$cell.ClearFormats()
$ws.Cells.Item(17, 5).Value = '  +2.23%  '

# Row 18
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.552.63'
$cell.ClearFormats()
$ws.Cells.Item(18, 5).Value = '  +1.64%  '

# Row 19
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = '11.04'
$cell.ClearFormats()
$ws.Cells.Item(19, 5).Value = '  -2.43%  '

# Row 20
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.27'
$cell.ClearFormats()
$ws.Cells.Item(20, 5).Value = '  -0.26%  '

# Row 21
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = '322.88'
$cell.ClearFormats()
$ws.Cells.Item(21, 5).Value = '  +0.39%  '

# Row 22
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.ClearFormats()
$ws.Cells.Item(22, 5).Value = '  +0.02%  '

# Row 23
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.86'
$cell.ClearFormats()
$ws.Cells.Item(23, 5).Value = '  +2.01%  '

# Row 24
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = '62.98'
$cell.ClearFormats()
$ws.Cells.Item(24, 5).Value = '  +3.76%  '

# Row 25
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.423'
$cell.ClearFormats()
$ws.Cells.Item(25, 5).Value = '  -3.21%  '

# Row 26
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.165'
$cell.ClearFormats()
$ws.Cells.Item(26, 5).Value = '  +2.38%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  +0.35%  '

# Row 28
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.88'
$cell.ClearFormats()
$ws.Cells.Item(28, 5).Value = '  +2.06%  '

# Row 29
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = '6.91'
$cell.ClearFormats()
$ws.Cells.Item(29, 5).Value = '  +0.56%  '

# Row 30
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0₃0776'
$cell.ClearFormats()
$ws.Cells.Item(30, 5).Value = '  +0.97%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  +0.92%  '

# Row 32
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = '164.36'
$cell.ClearFormats()
$ws.Cells.Item(32, 5).Value = '  +4.97%  '

# Row 33
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.15'
$cell.ClearFormats()
$ws.Cells.Item(33, 5).Value = '  -8.18%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  -0.02%  '

# Row 35
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.44'
$cell.ClearFormats()
$ws.Cells.Item(35, 5).Value = '  +2.75%  '

# Row 36
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = '18.56'
$cell.ClearFormats()
$ws.Cells.Item(36, 5).Value = '  +0.45%  '

# Row 37
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.31'
$cell.ClearFormats()
$ws.Cells.Item(37, 5).Value = '  -1.53%  '

# Row 38
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.59'
$cell.ClearFormats()
$ws.Cells.Item(38, 5).Value = '  -0.54%  '

# Row 39
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = '36.97'
$cell.ClearFormats()
$ws.Cells.Item(39, 5).Value = '  +0.49%  '

# Row 40
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.47'
$cell.ClearFormats()
$ws.Cells.Item(40, 5).Value = '  -6.93%  '

# Row 41
$cell = $ws.Cells.Item(41, 2)
$cell.NumberFormat = "@"
$cell.Value = 'Filecoin'
$cell.ClearFormats()
$cell = $ws.Cells.Item(41, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$cell.ClearFormats()
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.67'
$cell.ClearFormats()
$ws.Cells.Item(41, 5).Value = '  +0.34%  '

# Row 42
$cell = $ws.Cells.Item(42, 2)
$cell.NumberFormat = "@"
$cell.Value = 'Bittensor'
$cell.ClearFormats()
$cell = $ws.Cells.Item(42, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$cell.ClearFormats()
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = '291.59'
$cell.ClearFormats()
$ws.Cells.Item(42, 5).Value = '  -5.58%  '

# Row 43
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.815'
$cell.ClearFormats()
$ws.Cells.Item(43, 5).Value = '  +1.63%  '

# Row 44
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.997'
$cell.ClearFormats()
$ws.Cells.Item(44, 5).Value = '  -0.07%  '

# Row 45
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.602'
$cell.ClearFormats()
$ws.Cells.Item(45, 5).Value = '  +1.24%  '

# Row 46
$ws.Cells.Item(46, 5).Value = '  +0.78%  '

# Row 47
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = '125.20'
$cell.ClearFormats()
$ws.Cells.Item(47, 5).Value = '  +0.52%  '

# Row 48
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0931'
$cell.ClearFormats()
$ws.Cells.Item(48, 5).Value = '  +0.80%  '

# Row 49
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = '18.73'
$cell.ClearFormats()
$ws.Cells.Item(49, 5).Value = '  +1.16%  '

# Row 50
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0512'
$cell.ClearFormats()
$ws.Cells.Item(50, 5).Value = '  -0.46%  '

# Row 51
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0225'
$cell.ClearFormats()
$ws.Cells.Item(51, 5).Value = '  -0.87%  '
